$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-146 down to 51-147.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new data point.
$ws.Cells.Item(50, 1).Value  = 10
$ws.Cells.Item(50, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(50, 3).Value  = "La Araucanía"
$ws.Cells.Item(50, 4).Value  = 44662
$ws.Cells.Item(50, 5).Value  = 9
$ws.Cells.Item(50, 6).Value  = "Fruta"
$ws.Cells.Item(50, 7).Value  = 100104
$ws.Cells.Item(50, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(50, 9).Value  = 100104003
$ws.Cells.Item(50, 10).Value = "Membrillo"
$ws.Cells.Item(50, 11).Value = "Champion"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 85
$ws.Cells.Item(50, 14).Value = 15000
$ws.Cells.Item(50, 15).Value = 15000
$ws.Cells.Item(50, 16).Value = 15000
$ws.Cells.Item(50, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 833
$ws.Cells.Item(50, 20).Value = 18
